$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.916650891304016
$ws.Range("B1").Value = 2.075867414474487
$ws.Range("C1").Value = 2.122825622558594
$ws.Range("D1").Value = 2.584294080734253
$ws.Range("E1").Value = 3.699689388275146
